$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header text in B1 from "Row" to "Column"
$ws.Range("B1").Value = "Column"

# Update J2:J12 values (column df7_bc) with newly computed figures
$ws.Range("J2").Value = 4.844860156132717
$ws.Range("J3").Value = 4.584225672555992
$ws.Range("J4").Value = 4.292236403526178
$ws.Range("J5").Value = 3.996000688478106
$ws.Range("J6").Value = 3.69343394498344
$ws.Range("J7").Value = 3.389067992410719
$ws.Range("J8").Value = 3.089356884241885
$ws.Range("J9").Value = 2.799802626185456
$ws.Range("J10").Value = 2.510702835947995
$ws.Range("J11").Value = 2.230745030236432
$ws.Range("J12").Value = 1.937362681470768
